$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 13.57145098701222
$ws.Range("D2").Value = 6.5697511566171
$ws.Range("E2").Value = 13.95778867565
$ws.Range("F2").Value = 34.20796466813351
$ws.Range("G2").Value = 46.28226871543323
$ws.Range("H2").Value = 18.34431723902467
$ws.Range("K2").Value = 18.82246140342266
$ws.Range("L2").Value = 9.32667458863575
$ws.Range("N2").Value = 18.19516037576441

$ws.Range("C3").Value = 13.48164347717125
$ws.Range("D3").Value = 6.596071908262899
$ws.Range("E3").Value = 13.90521294087131
$ws.Range("F3").Value = 34.1219364324819
$ws.Range("G3").Value = 46.08486193093518
$ws.Range("H3").Value = 18.37656594735781
$ws.Range("K3").Value = 18.38193157570318
$ws.Range("L3").Value = 9.332968982502484
$ws.Range("N3").Value = 18.26770105833913

$ws.Range("C4").Value = 13.42983384161126
$ws.Range("D4").Value = 6.612877010537989
$ws.Range("E4").Value = 13.87627394787984
$ws.Range("F4").Value = 34.08032285050331
$ws.Range("G4").Value = 45.9808008251845
$ws.Range("H4").Value = 18.40086016960236
$ws.Range("K4").Value = 18.11034350240975
$ws.Range("L4").Value = 9.338469991508651
$ws.Range("N4").Value = 18.31423965276836

$ws.Range("C5").Value = 13.40957632625182
$ws.Range("D5").Value = 6.619887987324041
$ws.Range("E5").Value = 13.86532889080523
$ws.Range("F5").Value = 34.06618738580693
$ws.Range("G5").Value = 45.94272713960191
$ws.Range("H5").Value = 18.41188570315701
$ws.Range("K5").Value = 17.99955941111204
$ws.Range("L5").Value = 9.341122885524577
$ws.Range("N5").Value = 18.33370890699295

$ws.Range("C6").Value = 13.40626471764774
$ws.Range("D6").Value = 6.621062008721415
$ws.Range("E6").Value = 13.86356286961531
$ws.Range("F6").Value = 34.06401071683383
$ws.Range("G6").Value = 45.93666712269791
$ws.Range("H6").Value = 18.41378432169954
$ws.Range("K6").Value = 17.98116215798089
$ws.Range("L6").Value = 9.341588221558379
$ws.Range("N6").Value = 18.33697228232026

$ws.Range("C7").Value = 13.429557156696
$ws.Range("D7").Value = 6.612970902973049
$ws.Range("E7").Value = 13.87612289743222
$ws.Range("F7").Value = 34.08012078437681
$ws.Range("G7").Value = 45.98026978890487
$ws.Range("H7").Value = 18.40100431317611
$ws.Range("K7").Value = 18.10884964902587
$ws.Range("L7").Value = 9.338504104909106
$ws.Range("N7").Value = 18.31450017716789

$ws.Range("C8").Value = 13.53980455102993
$ws.Range("D8").Value = 6.578693367682603
$ws.Range("E8").Value = 13.93897145824119
$ws.Range("F8").Value = 34.17597938630824
$ws.Range("G8").Value = 46.21065826766266
$ws.Range("H8").Value = 18.354501405858
$ws.Range("K8").Value = 18.67089369504777
$ws.Range("L8").Value = 9.328505176571888
$ws.Range("N8").Value = 18.21975880111536

$ws.Range("C9").Value = 13.78160960751934
$ws.Range("D9").Value = 6.516548030361411
$ws.Range("E9").Value = 14.08837333526263
$ws.Range("F9").Value = 34.45254881122508
$ws.Range("G9").Value = 46.79726910751912
$ws.Range("H9").Value = 18.29914866099062
$ws.Range("K9").Value = 19.75729905626624
$ws.Range("L9").Value = 9.321887760973713
$ws.Range("N9").Value = 18.04973931025036

$ws.Range("C10").Value = 13.97367028911398
$ws.Range("D10").Value = 6.473929052749505
$ws.Range("E10").Value = 14.21352606907382
$ws.Range("F10").Value = 34.70906210091638
$ws.Range("G10").Value = 47.30829347761023
$ws.Range("H10").Value = 18.28056482489473
$ws.Range("K10").Value = 20.53699256931422
$ws.Range("L10").Value = 9.324950504643704
$ws.Range("N10").Value = 17.93431464791041

$ws.Range("C11").Value = 14.06389087340775
$ws.Range("D11").Value = 6.455189010960082
$ws.Range("E11").Value = 14.27366324696787
$ws.Range("F11").Value = 34.83711304726128
$ws.Range("G11").Value = 47.55755003576834
$ws.Range("H11").Value = 18.27694852444758
$ws.Range("K11").Value = 20.8859726428966
$ws.Range("L11").Value = 9.328063150234385
$ws.Range("N11").Value = 17.88383892282582

$ws.Range("C12").Value = 14.09844001095476
$ws.Range("D12").Value = 6.448184882973342
$ws.Range("E12").Value = 14.29688327503047
$ws.Range("F12").Value = 34.88721232519542
$ws.Range("G12").Value = 47.65429088724133
$ws.Range("H12").Value = 18.27627754819133
$ws.Range("K12").Value = 21.01717204442344
$ws.Range("L12").Value = 9.329488684842516
$ws.Range("N12").Value = 17.86501523065189

$ws.Range("C13").Value = 14.09098254532692
$ws.Range("D13").Value = 6.44968925401925
$ws.Range("E13").Value = 14.29186275150074
$ws.Range("F13").Value = 34.87635140838341
$ws.Range("G13").Value = 47.63335241176232
$ws.Range("H13").Value = 18.2763909516772
$ws.Range("K13").Value = 20.9889602221643
$ws.Range("L13").Value = 9.329170699926159
$ws.Range("N13").Value = 17.86905636318856

$ws.Range("C14").Value = 14.06672569287845
$ws.Range("D14").Value = 6.45461093148967
$ws.Range("E14").Value = 14.27556468623266
$ws.Range("F14").Value = 34.84120264775933
$ws.Range("G14").Value = 47.56546218258472
$ws.Range("H14").Value = 18.2768793105821
$ws.Range("K14").Value = 20.8967862886527
$ws.Range("L14").Value = 9.328175485517102
$ws.Range("N14").Value = 17.88228447694939

$ws.Range("C15").Value = 14.05191698393129
$ws.Range("D15").Value = 6.457637600066104
$ws.Range("E15").Value = 14.26563950070359
$ws.Range("F15").Value = 34.81988179023484
$ws.Range("G15").Value = 47.52418193905644
$ws.Range("H15").Value = 18.27726947753769
$ws.Range("K15").Value = 20.84019941515877
$ws.Range("L15").Value = 9.327598019918202
$ws.Range("N15").Value = 17.8904248412021

$ws.Range("C16").Value = 13.96782937796532
$ws.Range("D16").Value = 6.475166720674573
$ws.Range("E16").Value = 14.20965931819968
$ws.Range("F16").Value = 34.70092026685781
$ws.Range("G16").Value = 47.29233669523259
$ws.Range("H16").Value = 18.28089875224891
$ws.Range("K16").Value = 20.51405963601406
$ws.Range("L16").Value = 9.324781653167454
$ws.Range("N16").Value = 17.93765408386598

$ws.Range("C17").Value = 13.91695676454551
$ws.Range("D17").Value = 6.486085553003484
$ws.Range("E17").Value = 14.17612876813031
$ws.Range("F17").Value = 34.63083476497768
$ws.Range("G17").Value = 47.1543653837258
$ws.Range("H17").Value = 18.28436634274221
$ws.Range("K17").Value = 20.31242546875508
$ws.Range("L17").Value = 9.323494037440746
$ws.Range("N17").Value = 17.96714675972008

$ws.Range("C18").Value = 13.88796612904851
$ws.Range("D18").Value = 6.492426774220193
$ws.Range("E18").Value = 14.1571453709264
$ws.Range("F18").Value = 34.59159455154342
$ws.Range("G18").Value = 47.07659226357895
$ws.Range("H18").Value = 18.28681600447562
$ws.Range("K18").Value = 20.19592094790791
$ws.Range("L18").Value = 9.322915304796954
$ws.Range("N18").Value = 17.98430148874661

$ws.Range("C19").Value = 13.87819752415979
$ws.Range("D19").Value = 6.494584303440516
$ws.Range("E19").Value = 14.15077026697977
$ws.Range("F19").Value = 34.578493150815
$ws.Range("G19").Value = 47.05053350660186
$ws.Range("H19").Value = 18.28772350569523
$ws.Range("K19").Value = 20.15638766800898
$ws.Range("L19").Value = 9.322747167156113
$ws.Range("N19").Value = 17.99014269887855

$ws.Range("C20").Value = 13.922344484597
$ws.Range("D20").Value = 6.484916917474068
$ws.Range("E20").Value = 14.17966694069282
$ws.Range("F20").Value = 34.63818479852981
$ws.Range("G20").Value = 47.16888912027321
$ws.Range("H20").Value = 18.28395007938213
$ws.Range("K20").Value = 20.33394559204799
$ws.Range("L20").Value = 9.323614355915687
$ws.Range("N20").Value = 17.96398742586264

$ws.Range("C21").Value = 14.07384028792216
$ws.Range("D21").Value = 6.453162816349183
$ws.Range("E21").Value = 14.28033979227506
$ws.Range("F21").Value = 34.85148323113272
$ws.Range("G21").Value = 47.58533985239232
$ws.Range("H21").Value = 18.27671689213394
$ws.Range("K21").Value = 20.92388683568551
$ws.Range("L21").Value = 9.328461108754876
$ws.Range("N21").Value = 17.87839119270723

$ws.Range("C22").Value = 14.17508078834965
$ws.Range("D22").Value = 6.432947412241244
$ws.Range("E22").Value = 14.34873622676577
$ws.Range("F22").Value = 35.00025182268543
$ws.Range("G22").Value = 47.8711959589408
$ws.Range("H22").Value = 18.27606161151895
$ws.Range("K22").Value = 21.30384584991977
$ws.Range("L22").Value = 9.333067248573771
$ws.Range("N22").Value = 17.82414079456741

$ws.Range("C23").Value = 14.12085143668849
$ws.Range("D23").Value = 6.443687814263292
$ws.Range("E23").Value = 14.31199842330867
$ws.Range("F23").Value = 34.92000323776512
$ws.Range("G23").Value = 47.71739907967425
$ws.Range("H23").Value = 18.27603794308786
$ws.Range("K23").Value = 21.10160796199213
$ws.Range("L23").Value = 9.330477410494204
$ws.Range("N23").Value = 17.85294104174629

$ws.Range("C24").Value = 13.91990789485794
$ws.Range("D24").Value = 6.485445058658898
$ws.Range("E24").Value = 14.1780664168623
$ws.Range("F24").Value = 34.63485856687741
$ws.Range("G24").Value = 47.16231811115887
$ws.Range("H24").Value = 18.28413685151591
$ws.Range("K24").Value = 20.32421814810556
$ws.Range("L24").Value = 9.323559456749159
$ws.Range("N24").Value = 17.96541514067841

$ws.Range("C25").Value = 13.71357282499553
$ws.Range("D25").Value = 6.532822424983929
$ws.Range("E25").Value = 14.04520963879768
$ws.Range("F25").Value = 34.36830218028785
$ws.Range("G25").Value = 46.62433590629362
$ws.Range("H25").Value = 18.31026042467342
$ws.Range("K25").Value = 19.46603407518812
$ws.Range("L25").Value = 9.322286172935479
$ws.Range("N25").Value = 18.0940588814531
